$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to B:F
$ws.Columns.Item(1).Insert()

# Populate the new column A with the staging URL header/value
$ws.Range("A1").Value = "ESRIActStgURL"
$ws.Range("A2").Value = "https://accounts-stg.esri.com/"

# Add a hyperlink on A2 pointing at the staging URL
$ws.Hyperlinks.Add($ws.Range("A2"), "https://accounts-stg.esri.com/")

# The Hyperlinks.Add call stamps a "Hyperlink" style on the cell; restore
# the default so A2 keeps the sheet's plain formatting
$ws.Range("A2").Style = "Normal"

# Match the saved selection from the edited workbook
$ws.Range("A1:A2").Select()
